# Update the "correlation" (column G) values as part of the
# "reflecting proper edits to headscan_full1" reknit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G7").Value = 0.9453
$ws.Range("G17").Value = 0.9453
$ws.Range("G28").Value = 0.737
$ws.Range("G33").Value = 0.8045
$ws.Range("G41").Value = 0.8082
$ws.Range("G43").Value = 0.8082
$ws.Range("G44").Value = 0.9453
$ws.Range("G45").Value = 0.8045
$ws.Range("G46").Value = 0.737
